$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Name = "LoginData"

# Header row
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Password"

# valid
$ws.Range("A2").Value = "valid"
$ws.Range("B2").Value = "testvaliduser@gmail.com"
$ws.Range("C2").Value = "Test@123"

# invalid
$ws.Range("A3").Value = "invalid"
$ws.Range("B3").Value = "wrong@gmail.com"
$ws.Range("C3").Value = "wrong123"

# blank
$ws.Range("A4").Value = "blank"

$ws.Range("A1:C3").Font.ThemeColor = 1
$ws.Range("A4").Font.ThemeColor = 1

$ws.Columns.Item(1).ColumnWidth = 25.83
$ws.Columns.Item(2).ColumnWidth = 28.66
$ws.Columns.Item(3).ColumnWidth = 20.16
